$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 26 (Email / varchar / 50 / NEIN) — content below
# shifts up by one row.
$ws.Rows("26:26").Delete()

# Row-delete shifts the row references used by the two far-below defined
# names down by one as well.
$n1 = $wb.Names.Item("A2000381")
$n1.RefersTo = "=Tabelle1!`$A`$1000378"

$n2 = $wb.Names.Item("A9999999")
$n2.RefersTo = "=Tabelle1!`$A`$99996"

# Leave the selection on the (new) row 26, matching Excel's behaviour
# after a row deletion.
$ws.Rows("26:26").Select()
